# 20220509 Docs : 첫 번째 ppt 업데이트
#
# 1) The "datetimeFigureOut" date placeholder (slide master + every slide
#    layout) auto-updates from 2022-05-08 to 2022-05-09.
# 2) Fix the typo "Envionment" -> "Environment" in the title of slide 5.

$p = $ppt.ActivePresentation

$oldDate = "2022-05-08"
$newDate = "2022-05-09"

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master date placeholder.
Update-DatePlaceholder -shapes $p.SlideMaster.Shapes

# Every slide layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DatePlaceholder -shapes $layouts.Item($L).Shapes
}

# Slide 5 title typo fix: "Envionment" -> "Environment".
$s5 = $p.Slides.Item(5)
for ($i = 1; $i -le $s5.Shapes.Count; $i++) {
    $shp = $s5.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "Envionment") {
        $shp.TextFrame.TextRange.Text = "Environment"
    }
}
